$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("C2").Value = $true
$ws.Range("D2").Value = 0.02130412373605337
$ws.Range("E2").Value = 0.02130412373605337

# Row 3
$ws.Range("D3").Value = [double]"5.637221135543975E-16"
$ws.Range("E3").Value = [double]"5.637221135543975E-16"

# Row 4
$ws.Range("D4").Value = [double]"6.136143485490833E-35"
$ws.Range("E4").Value = [double]"6.136143485490833E-35"

# Row 5
$ws.Range("D5").Value = [double]"1.120076490148989E-53"
$ws.Range("E5").Value = [double]"1.120076490148989E-53"

# Row 6
$ws.Range("D6").Value = [double]"1.385832130780251E-09"
$ws.Range("E6").Value = [double]"1.385832130780251E-09"

# Row 7
$ws.Range("D7").Value = [double]"1.224215514238188E-24"

# Row 9
$ws.Range("D9").Value = 0.849178715801851
$ws.Range("E9").Value = 0.150821284198149

# Row 10
$ws.Range("D10").Value = [double]"5.870366897997246E-10"
$ws.Range("E10").Value = 0.9999999994129634

# Row 11
$ws.Range("D11").Value = [double]"5.313198159106881E-08"
$ws.Range("E11").Value = 0.9999999468680184
$ws.Range("F11").Value = 9.325117111206055
$ws.Range("G11").Value = 0.7
